# ADDED: Negative Test Cases
#
# Duplicate the "BooksWagon" sheet into a new "invalid" sheet (used to hold
# a negative/invalid test row), trim it down to a single data row with
# different sample values, and touch up the original sheet's B3 cell +
# view/selection state.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Copy BooksWagon -> new sheet placed right after it, then rename it.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "invalid"

# 2) Drop the extra data row (row 3) on the new sheet - only one sample
#    row is kept there. Hyperlinks on this engine are sheet-scoped when
#    deleted via a Range, so clear them all and re-add only the ones that
#    remain (row 2).
$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Rows.Item(3).Delete()

$ws2.Range("A2").Value = "dineshkumar.icon@gmail.com"
$ws2.Range("B2").Value = "Dinnu@248"
$ws2.Range("D2").Value = -1

$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:dineshkumar.icon@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Dinnu@248") | Out-Null
$ws2.Range("A2:B2").Style = "Hyperlink"

# 3) New sheet's column widths + view/selection.
$ws2.Columns.Item(1).ColumnWidth = 47.42578125
$ws2.Columns.Item(2).ColumnWidth = 37.42578125
$ws2.Columns.Item(3).ColumnWidth = 32.140625
$ws2.Columns.Item(4).ColumnWidth = 28.28515625

$ws2.Range("D3").Select()

# 4) Back on the original sheet: B3 becomes "Dinnu@247" (was the stray
#    "`" value), the matching hyperlink no longer needs an explicit
#    display override, and the view scrolls right with a block selection.
$ws1.Range("B3").Value = "Dinnu@247"
$ws1.Hyperlinks.Item(4).TextToDisplay = "Dinnu@247"

$ws1.Range("A1:L2").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 9
